# Generate Report for Handback
#
# This script updates the localization-status report to reflect a failed
# handback transform:
#   - The "Ready for handoff" status (shown on the Overview sheet and on
#     each language sheet for the 13089187-... file) becomes
#     "Handback transform failed".
#   - Each language sheet gets an "Error Detail" entry (column K, row 3)
#     explaining the handback file name mismatch.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: zhjrlhkq.4pt is different with handoff file name: 13089187-b6dd-400e-95f6-1ead192e4cab.7a67af721f4bdfdbcbc08e97ac7f76a184dfb813.zh-cn."

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: zhjrlhkq.4pt is different with handoff file name: 13089187-b6dd-400e-95f6-1ead192e4cab.7a67af721f4bdfdbcbc08e97ac7f76a184dfb813.de-de."
